# B6-PowerPoint.pptx — Sat, Jun 06, 2020 6:04:53 PM
#
# 1) Swap the deck's applied colour scheme from the "Integral / Red Violet"
#    palette over to the stock "Office" palette (the theme part backing the
#    slide master, theme1.xml, is the only theme surface this host's COM
#    layer can address from a slide-anchored object — see notes below).
# 2) Re-style the three tables (slides 14-16) from the custom "Table_0"
#    style onto the built-in table style {4A1AC524-B407-4667-9A78-312A06FE66E7}.

$p = $ppt.ActivePresentation

# --- 1) Theme colours -------------------------------------------------
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink. RGB is packed VBA-style as 0xBBGGRR, so a target hex "RRGGBB"
# becomes 0x(BB)(GG)(RR).
$officeColors = @(
    0x000000,  #  1 dk1      000000
    0xFFFFFF,  #  2 lt1      FFFFFF
    0x6A5444,  #  3 dk2      44546A
    0xE6E6E7,  #  4 lt2      E7E6E6
    0xD59B5B,  #  5 accent1  5B9BD5
    0x317DED,  #  6 accent2  ED7D31
    0xA5A5A5,  #  7 accent3  A5A5A5
    0x00C0FF,  #  8 accent4  FFC000
    0xC47244,  #  9 accent5  4472C4
    0x47AD70,  # 10 accent6  70AD47
    0xC16305,  # 11 hlink    0563C1
    0x724F95   # 12 folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}

# --- 2) Table styles ----------------------------------------------------
$newTableStyle = "{4A1AC524-B407-4667-9A78-312A06FE66E7}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $s.Shapes.Count; $shapeIdx++) {
        $sh = $s.Shapes.Item($shapeIdx)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyle)
        }
    }
}
